# New test case that shows that words which style is changed at midword
# are returned cut in two parts.
#
# Sheet1 ("First sheet") column D currently holds the plain string
# "spreadsheet". We re-type it and then bold just the second half
# ("sheet") using Range.Characters(start, length), which is exactly how
# Excel lets a user apply character-level (mid-word) formatting to part
# of a cell's text. We also add a brand new cell (E1) containing another
# word ("midword") that is half bold in the same way, to exercise the
# same mid-word-split behaviour on a second, brand-new shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First sheet")

# --- D1: "spreadsheet" -> "spread" + bold "sheet" ---------------------
$d1 = $ws.Range("D1")
$d1.Value = "spreadsheet"
$chars = $d1.Characters(7, 5)
$chars.Font.Name = "Arial"
$chars.Font.Size = 10
$chars.Font.Bold = $true

# --- E1 (new cell): "midword" -> "mid" + bold "word" -------------------
$e1 = $ws.Range("E1")
$e1.Value = "midword"
$chars2 = $e1.Characters(4, 4)
$chars2.Font.Name = "Arial"
$chars2.Font.Size = 10
$chars2.Font.Bold = $true
